$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: append a brand-new, empty paragraph at the very end of the document
# and return it (as the new $d.Paragraphs.Last).
# ---------------------------------------------------------------------------
function New-EndParagraph {
    $lp = $d.Paragraphs.Last
    $rr = $lp.Range
    $rr.Collapse(0)
    $rr.InsertParagraphAfter()
    return $d.Paragraphs.Last
}

# A paragraph that already carries the shared bullet numbering definition
# (abstractNumId 991) - used as a template so the new bulleted lists reuse
# the same abstract numbering definition instead of minting a new one.
$srcListPara = $d.Paragraphs.Item(66)

# =====================================================================
# 17 Team Members - Physicians
# =====================================================================
$p = New-EndParagraph
$p.Range.Text = "17 Team Members - Physicians"
$p.Style = "Heading2"
$sec1First = $p

$p = New-EndParagraph
$p.Range.Text = "Primary Care Provider"
$p.Style = "FirstParagraph"

$p = New-EndParagraph
$p.Range.Text = "Gastroenterologist"
$p.Style = "BodyText"

$p = New-EndParagraph
$p.Range.Text = "Medical Oncologist (chemotherapy)"
$p.Style = "BodyText"

$p = New-EndParagraph
$p.Range.Text = "Radiation Oncologist (radiation)"
$p.Style = "BodyText"

$p = New-EndParagraph
$p.Range.Text = "Surgeons"
$p.Style = "BodyText"

# Bulleted sub-list of surgeons - first item mints a fresh numId (1012) that
# points at the same abstract numbering definition the rest of the doc uses.
$tmpl1 = $srcListPara.Range.ListFormat.ListTemplate
$p = New-EndParagraph
$p.Range.Text = "Jonathan Salo MD"
$p.Style = "Compact"
$p.Range.ListFormat.ApplyListTemplateWithLevel($tmpl1)

# Subsequent bullets: insert directly after the previous bullet's range so
# they inherit the very same numId instead of minting new ones.
$rPrev = $p.Range
$rPrev.Collapse(0)
$rPrev.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Jeffrey Hagen MD"

$rPrev = $p.Range
$rPrev.Collapse(0)
$rPrev.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Michael Roach MD"
$sec1Last = $p

# Bookmark wrapping the whole "Team Members - Physicians" section.
$bm1Range = $d.Range($sec1First.Range.Start, $sec1Last.Range.End)
$d.Bookmarks.Add("team-members---physicians", $bm1Range)

# =====================================================================
# 18 Team Members - Support Staff
# =====================================================================
$p = New-EndParagraph
$p.Range.Text = "18 Team Members - Support Staff"
$p.Style = "Heading2"
$sec2First = $p

$p = New-EndParagraph
$p.Range.Text = "Dietitian - Liz Koch"
$p.Style = "FirstParagraph"

$p = New-EndParagraph
$p.Range.Text = "Nurses"
$p.Style = "BodyText"

# Bulleted sub-list of nurses - mints another fresh numId (1013).
$tmpl2 = $srcListPara.Range.ListFormat.ListTemplate
$p = New-EndParagraph
$p.Range.Text = "Matthew Carpenter RN"
$p.Style = "Compact"
$p.Range.ListFormat.ApplyListTemplateWithLevel($tmpl2)

$rPrev = $p.Range
$rPrev.Collapse(0)
$rPrev.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Brandon Galloway LPN"

$p = New-EndParagraph
$p.Range.Text = "Navigator - Laura Swift"
$p.Style = "FirstParagraph"
$sec2Last = $p

# Bookmark wrapping the whole "Team Members - Support Staff" section.
$bm2Range = $d.Range($sec2First.Range.Start, $sec2Last.Range.End)
$d.Bookmarks.Add("team-members---support-staff", $bm2Range)

Write-Output "Added sections 17 and 18 with bookmarks and bullet lists."
